$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("SFIA Level"), which pushes the
# existing SFIA Level / Keycode / Description columns one slot to the right.
$ws.Columns("B:B").Insert()

# Header for the newly inserted column
$ws.Range("B1").Value = "Skill Description"

# Fill in the full skill name for each SkillCode group (column A values)
$ws.Range("B2:B5").Value = "Autonomy"
$ws.Range("B6:B8").Value = "Influence"
$ws.Range("B9:B11").Value = "Complexity"
$ws.Range("B12:B14").Value = "Knowledge"
$ws.Range("B15:B21").Value = "Information security"
$ws.Range("B22:B27").Value = "Information assurance"
$ws.Range("B28:B29").Value = "Risk management"
$ws.Range("B31:B34").Value = "Measurement"
